# Added the ability to attach VMM domain to EPGs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TenantConfig")

# Make the header row taller to fit the wrapped header text.
$ws.Rows.Item(1).RowHeight = 32

# Header E1: clarify that this column also accepts a VPC name, and wrap it.
$ws.Range("E1").Value = "Static_Path 1/101/1/13 or VPC name"
$ws.Range("E1").WrapText = $true

# Sample EPG row: new physical-domain value + a VMM domain example, and a
# tweaked encap_vlan sample value.
$ws.Range("F2").Value = 740
$ws.Range("H2").Value = "e7vmw1_Data"
$ws.Range("G2").Value = "e7_f5_phys"

# The old VMM_Domain sample value on row 3 is no longer needed.
$ws.Range("H3").ClearContents()

# Select TenantConfig / F2, matching the new active view.
$ws.Activate()
$ws.Range("F2").Select()
